$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '无'
$ws.Range("B2").Value = '无'
$ws.Range("C2").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤B3887学.jpg'
$ws.Range("A3").Value = '粤C1017学'
$ws.Range("B3").Value = 'yellow'
$ws.Range("C3").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1017学.jpg'
$ws.Range("A4").Value = '无'
$ws.Range("B4").Value = '无'
$ws.Range("C4").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1035学.jpg'
$ws.Range("A5").Value = '粤C1067学'
$ws.Range("B5").Value = 'yellow'
$ws.Range("C5").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1067学.jpg'
$ws.Range("A6").Value = '黑C1142学'
$ws.Range("B6").Value = 'yellow'
$ws.Range("C6").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1142学.jpg'
$ws.Range("A7").Value = '辽C1315学'
$ws.Range("B7").Value = 'yellow'
$ws.Range("C7").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1315学.jpg'
$ws.Range("A8").Value = '无'
$ws.Range("B8").Value = '无'
$ws.Range("C8").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1319学.jpg'
$ws.Range("A9").Value = '粤C1332学'
$ws.Range("B9").Value = 'yellow'
$ws.Range("C9").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1332学.jpg'
$ws.Range("A10").Value = '甘C1447学'
$ws.Range("B10").Value = 'yellow'
$ws.Range("C10").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1447学.jpg'
$ws.Range("A11").Value = '粤C1515学'
$ws.Range("B11").Value = 'yellow'
$ws.Range("C11").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1515学.jpg'
$ws.Range("A12").Value = '粤C1557学'
$ws.Range("B12").Value = 'yellow'
$ws.Range("C12").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1557学.jpg'
$ws.Range("A13").Value = '粤C1616学'
$ws.Range("B13").Value = 'yellow'
$ws.Range("C13").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1616学.jpg'
$ws.Range("A14").Value = '无'
$ws.Range("B14").Value = '无'
$ws.Range("C14").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1655学.jpg'
$ws.Range("A15").Value = '无'
$ws.Range("B15").Value = '无'
$ws.Range("C15").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1833学.jpg'
$ws.Range("A16").Value = '粤C1818辽'
$ws.Range("B16").Value = 'yellow'
$ws.Range("C16").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1881学.jpg'
$ws.Range("A17").Value = '黑C1889学'
$ws.Range("B17").Value = 'yellow'
$ws.Range("C17").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1889学.jpg'
$ws.Range("A18").Value = '粤C1926学'
$ws.Range("B18").Value = 'yellow'
$ws.Range("C18").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C1926学.jpg'
$ws.Range("A19").Value = '粤C2001学1'
$ws.Range("B19").Value = 'yellow'
$ws.Range("C19").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2001学.jpg'
$ws.Range("A20").Value = '粤C202J学'
$ws.Range("B20").Value = 'yellow'
$ws.Range("C20").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2021学.jpg'
$ws.Range("A21").Value = '粤C2027学'
$ws.Range("B21").Value = 'yellow'
$ws.Range("C21").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2027学.jpg'
$ws.Range("A22").Value = '粤C2028学'
$ws.Range("B22").Value = 'yellow'
$ws.Range("C22").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2028学.jpg'
$ws.Range("A23").Value = '甘C2043学'
$ws.Range("B23").Value = 'yellow'
$ws.Range("C23").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2043学.jpg'
$ws.Range("A24").Value = '粤C2062学'
$ws.Range("B24").Value = 'yellow'
$ws.Range("C24").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2062学.jpg'
$ws.Range("A25").Value = '粤C2073学'
$ws.Range("B25").Value = 'yellow'
$ws.Range("C25").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2073学.jpg'
$ws.Range("A26").Value = '甘C2078学'
$ws.Range("B26").Value = 'yellow'
$ws.Range("C26").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2078学.jpg'
$ws.Range("A27").Value = '无'
$ws.Range("B27").Value = '无'
$ws.Range("C27").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2091学.jpg'
$ws.Range("A28").Value = '粤C2092学'
$ws.Range("B28").Value = 'yellow'
$ws.Range("C28").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2092学.jpg'
$ws.Range("A29").Value = '无'
$ws.Range("B29").Value = '无'
$ws.Range("C29").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2110学.jpg'
$ws.Range("A30").Value = '警C2128学'
$ws.Range("B30").Value = 'yellow'
$ws.Range("C30").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2128学.jpg'

# Rows 31-51: clear A and B cells entirely, set C
$ws.Range("A31:B51").ClearContents()
$ws.Range("C31").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2219学.jpg'
$ws.Range("C32").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2231学.jpg'
$ws.Range("C33").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2239学.jpg'
$ws.Range("C34").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2249学.jpg'
$ws.Range("C35").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2273学.jpg'
$ws.Range("C36").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2322学.jpg'
$ws.Range("C37").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2339学.jpg'
$ws.Range("C38").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2372学.jpg'
$ws.Range("C39").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2402学.jpg'
$ws.Range("C40").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2415学.jpg'
$ws.Range("C41").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2427学.jpg'
$ws.Range("C42").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2523学.jpg'
$ws.Range("C43").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2542学.jpg'
$ws.Range("C44").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2557学.jpg'
$ws.Range("C45").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2559学.jpg'
$ws.Range("C46").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2586学.jpg'
$ws.Range("C47").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2618学.jpg'
$ws.Range("C48").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2661学.jpg'
$ws.Range("C49").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2687学.jpg'
$ws.Range("C50").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2725学.jpg'
$ws.Range("C51").Value = 'D:/third_party/openalpr-2.3.0/Task3_车牌识别/功能评测图像库/车牌种类变化子库/教练车牌/粤C2738学.jpg'
